# Weekly update: insert a new price record as row 22, pushing the
# existing rows 22-45 down to 23-46 (dimension grows from R45 to R46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 22 (this copies the
# formatting of row 22, including the date-number-format style on D,
# and shifts all rows 22..45 down to 23..46).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Cells.Item(22, 1).Value  = 1
$ws.Cells.Item(22, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value  = 44477
$ws.Cells.Item(22, 5).Value  = 15
$ws.Cells.Item(22, 6).Value  = 100112021
$ws.Cells.Item(22, 7).Value  = "Ají"
$ws.Cells.Item(22, 8).Value  = "Inferno"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 120
$ws.Cells.Item(22, 11).Value = 36000
$ws.Cells.Item(22, 12).Value = 37000
$ws.Cells.Item(22, 13).Value = 36500
$ws.Cells.Item(22, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 2433
$ws.Cells.Item(22, 17).Value = 15
$ws.Cells.Item(22, 18).Value = "Hortaliza"
